# Auto-generated edit script: bulk update of market-price-derived columns (H:N)
# across multiple worksheets, matching a scheduled market-data refresh run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3801.6667
$ws.Range("I62").Value = 3562
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 3562
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -2938
$ws.Range("N62").Value = -6248
$ws.Range("H65").Value = 3801.6667
$ws.Range("I65").Value = 3562
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 17810
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -14690
$ws.Range("N65").Value = -31240
$ws.Range("H92").Value = 1022.64
$ws.Range("I92").Value = 897.55554
$ws.Range("K92").Value = 897.55554
$ws.Range("M92").Value = 350.44446
$ws.Range("H107").Value = 437
$ws.Range("I107").Value = 388.84616
$ws.Range("J107").Value = 750
$ws.Range("K107").Value = 388.84616
$ws.Range("L107").Value = 750
$ws.Range("M107").Value = 1531.15384
$ws.Range("N107").Value = -4590
$ws.Range("H113").Value = 3004.6155
$ws.Range("I113").Value = 2074.3333
$ws.Range("J113").Value = 4273.1816
$ws.Range("K113").Value = 2074.3333
$ws.Range("L113").Value = 4273.1816
$ws.Range("M113").Value = 1179.6667
$ws.Range("N113").Value = -10781.1816
$ws.Range("H118").Value = 291.9
$ws.Range("I118").Value = 291.9
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 875.6999999999999
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = 781.3000000000001
$ws.Range("N118").ClearContents()
$ws.Range("H135").Value = 626.6842
$ws.Range("I135").Value = 384
$ws.Range("J135").Value = 2689.5
$ws.Range("K135").Value = 3456
$ws.Range("L135").Value = 24205.5
$ws.Range("M135").Value = -921
$ws.Range("N135").Value = -29275.5
$ws.Range("H137").Value = 914.21
$ws.Range("I137").Value = 769.7317
$ws.Range("J137").Value = 1014.61017
$ws.Range("K137").Value = 2309.1951
$ws.Range("L137").Value = 3043.83051
$ws.Range("M137").Value = 240.8049000000001
$ws.Range("N137").Value = -8143.83051
$ws.Range("H138").Value = 3639399.2
$ws.Range("I138").Value = 1562.2941
$ws.Range("J138").Value = 9529231
$ws.Range("K138").Value = 4686.8823
$ws.Range("L138").Value = 28587693
$ws.Range("M138").Value = 453.1176999999998
$ws.Range("N138").Value = -28597973

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1523.1482
$ws.Range("I110").Value = 1496.6818
$ws.Range("J110").Value = 1639.6
$ws.Range("K110").Value = 1496.6818
$ws.Range("L110").Value = 1639.6
$ws.Range("M110").Value = 548.3181999999999
$ws.Range("N110").Value = -5729.6
$ws.Range("H122").Value = 1349.6333
$ws.Range("I122").Value = 876.55554
$ws.Range("J122").Value = 2059.25
$ws.Range("K122").Value = 2629.66662
$ws.Range("L122").Value = 6177.75
$ws.Range("M122").Value = -179.66662
$ws.Range("N122").Value = -11077.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1921.8667
$ws.Range("I107").Value = 1975.2307
$ws.Range("J107").Value = 1575
$ws.Range("K107").Value = 1975.2307
$ws.Range("L107").Value = 1575
$ws.Range("M107").Value = -55.23070000000007
$ws.Range("N107").Value = -5415
$ws.Range("H134").Value = 10278.65
$ws.Range("I134").Value = 10871.6
$ws.Range("J134").Value = 8499.799999999999
$ws.Range("K134").Value = 32614.8
$ws.Range("L134").Value = 25499.4
$ws.Range("M134").Value = -30079.8
$ws.Range("N134").Value = -30569.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1805.5555
$ws.Range("I99").Value = 1650
$ws.Range("K99").Value = 1650
$ws.Range("M99").Value = -152
$ws.Range("H122").Value = 1569
$ws.Range("I122").Value = 1100
$ws.Range("J122").Value = 2507
$ws.Range("K122").Value = 3300
$ws.Range("L122").Value = 7521
$ws.Range("M122").Value = -850
$ws.Range("N122").Value = -12421
$ws.Range("H126").Value = 1805.5555
$ws.Range("I126").Value = 1650
$ws.Range("K126").Value = 4950
$ws.Range("M126").Value = -2480

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 842.6087
$ws.Range("I131").Value = 329
$ws.Range("J131").Value = 1313.4166
$ws.Range("K131").Value = 987
$ws.Range("L131").Value = 3940.2498
$ws.Range("M131").Value = 4053
$ws.Range("N131").Value = -14020.2498

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 657.7273
$ws.Range("I107").Value = 476.0909
$ws.Range("J107").Value = 839.36365
$ws.Range("K107").Value = 476.0909
$ws.Range("L107").Value = 839.36365
$ws.Range("M107").Value = 1443.9091
$ws.Range("N107").Value = -4679.36365
$ws.Range("H122").Value = 202160
$ws.Range("I122").Value = 336266.66
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 1008799.98
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -1006349.98
$ws.Range("N122").Value = -7900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H21").Value = 70007
$ws.Range("J21").Value = 70007
$ws.Range("L21").Value = 70007
$ws.Range("N21").Value = -70355
$ws.Range("H46").Value = 926.0741
$ws.Range("I46").Value = 1183.6666
$ws.Range("J46").Value = 852.4761999999999
$ws.Range("K46").Value = 1183.6666
$ws.Range("L46").Value = 852.4761999999999
$ws.Range("M46").Value = -995.6666
$ws.Range("N46").Value = -1228.4762
$ws.Range("H61").Value = 2176.5625
$ws.Range("I61").Value = 1779.9
$ws.Range("J61").Value = 2837.6667
$ws.Range("K61").Value = 1779.9
$ws.Range("L61").Value = 2837.6667
$ws.Range("M61").Value = -1577.9
$ws.Range("N61").Value = -3241.6667
$ws.Range("H93").Value = 65419.09
$ws.Range("J93").Value = 78544.89
$ws.Range("L93").Value = 78544.89
$ws.Range("N93").Value = -81040.89
$ws.Range("H113").Value = 2176.5625
$ws.Range("I113").Value = 1779.9
$ws.Range("J113").Value = 2837.6667
$ws.Range("K113").Value = 1779.9
$ws.Range("L113").Value = 2837.6667
$ws.Range("M113").Value = 390.0999999999999
$ws.Range("N113").Value = -7177.6667
$ws.Range("H136").Value = 1487.2181
$ws.Range("I136").Value = 1480.7805
$ws.Range("J136").Value = 1506.0714
$ws.Range("K136").Value = 4442.3415
$ws.Range("L136").Value = 4518.2142
$ws.Range("M136").Value = -1892.3415
$ws.Range("N136").Value = -9618.2142

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 17434.166
$ws.Range("I107").Value = 20720.4
$ws.Range("K107").Value = 62161.2
$ws.Range("M107").Value = -60241.2
$ws.Range("H132").Value = 4493.8965
$ws.Range("I132").Value = 5424.1763
$ws.Range("J132").Value = 3176
$ws.Range("K132").Value = 16272.5289
$ws.Range("L132").Value = 9528
$ws.Range("M132").Value = -13742.5289
$ws.Range("N132").Value = -14588
